$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.229.30"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "'1.588.04"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'211.82"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").Value = "'19.22"
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "'1.810.71"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "'1.590.18"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").Value = "'63.92"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "'26.230.52"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D19").Value = "'7.37"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("D20").Value = "'214.16"
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("D25").Value = "'144.72"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("D33").Value = "'1.424.25"
$ws.Range("E33").Value = "  +8.80%  "
$ws.Range("D34").Value = "'2.95"
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.588"
$ws.Range("E36").Value = "  -4.27%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'1.45"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("E39").Value = "  +5.12%  "
$ws.Range("D40").Value = "'0.823"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").Value = "  -13.44%  "
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "'1.722.33"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").Value = "'61.13"
$ws.Range("D47").Value = "'85.75"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("D49").Value = "'0.0501"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("E50").Value = "  -1.43%  "
